$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Account" row to the "User" sheet (foreign key to Account) ---
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("A7").Value = "username"
$wsUser.Range("B7").Value = "varchar(64)"
$wsUser.Range("C7").Value = "foreignKey"
$wsUser.Range("D7").Value = "Account(username)"
$wsUser.Range("D12").Select() | Out-Null

# --- 2. Add a new "Account" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAccount = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAccount.Name = "Account"

$wsAccount.Range("A1").Value = "username"
$wsAccount.Range("B1").Value = "varchar(64)"
$wsAccount.Range("C1").Value = "primaryKey"

$wsAccount.Range("A2").Value = "password"
$wsAccount.Range("B2").Value = "varchar(64)"

$wsAccount.Range("A3").Value = "admin"
$wsAccount.Range("B3").Value = "bool"

$wsAccount.Range("A1").Font.Underline = $true
$wsAccount.Range("H9").Select() | Out-Null
